# Maersk_Example_Dana.xlsx — "Deleted Demand for Hydrogen"
#
# The Hydrogen_Kasso demand series (column B on the Demand sheet) is
# zeroed out, and the values move into the E-Methanol_Kasso demand
# series (column C), which was previously all zero. The "t" (time)
# column label also moves from the parameter-name row (row 3) up into
# the header row (row 1).
#
# The workbook's active tab is also moved back to the first sheet
# ("Units").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# "t" label moves from A3 up to A1 (header row, above the object names).
$ws.Range("A1").Value = "t"
$ws.Range("A3").Value = ""

# Hydrogen_Kasso (col B) demand deleted -> zeroed; the same values are
# now used for E-Methanol_Kasso (col C) demand, which used to be 0.
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 20

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 30

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 20

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 25

# Restore focus to the first sheet ("Units") as the active tab.
$wb.Worksheets.Item("Units").Activate()
